# Update PSSM numeric matrix (B2:K21) with supplemental-figure values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, -18.92234920697049, 1.941296412328727, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(3, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(4, -18.92234920697049, 1.980309275464383, 1.679498641772591, -18.92234920697049, 3.431780969387169, -18.92234920697049, 1.29510802035782, -18.92234920697049, 1.178632597196182, -18.92234920697049),
    @(5, -18.92234920697049, 1.703841963745995, -18.92234920697049, -18.92234920697049, -18.92234920697049, 2.927918910032342, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(6, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(7, 2.448949401647102, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(8, -18.92234920697049, -18.92234920697049, -18.92234920697049, 1.828103454887029, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(9, 3.861933574835612, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(10, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, 4.32192533619127, -18.92234920697049, 2.177720843817823),
    @(11, -18.92234920697049, -18.92234920697049, -18.92234920697049, 2.904343962751392, -18.92234920697049, 2.745305077094949, -18.92234920697049, -18.92234920697049, -18.92234920697049, 1.963044061103659),
    @(12, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(13, -18.92234920697049, -18.92234920697049, -18.92234920697049, 2.53287214261344, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, 1.82246574679058, 1.608423831000558),
    @(14, -18.92234920697049, -18.92234920697049, 1.500139599076208, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, 2.019068846174621),
    @(15, -18.92234920697049, -18.92234920697049, 1.706282332744112, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(16, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, 2.009494169584412, -18.92234920697049),
    @(17, -18.92234920697049, 2.146959599823662, 1.864906938847582, -18.92234920697049, -18.92234920697049, -18.92234920697049, 2.167634598019158, -18.92234920697049, 2.392340791972641, -18.92234920697049),
    @(18, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, -18.92234920697049, 2.033653232654954, -18.92234920697049, 2.299538184918097, -18.92234920697049),
    @(19, -18.92234920697049, -18.92234920697049, 2.085601273530149, -18.92234920697049, -18.92234920697049, -18.92234920697049, 1.66180314269102, -18.92234920697049, -18.92234920697049, -18.92234920697049),
    @(20, -18.92234920697049, 1.021460487237059, 1.494647220020012, -18.92234920697049, 3.203010341556001, -18.92234920697049, 1.775196184576318, -18.92234920697049, -18.92234920697049, 2.161676039943308),
    @(21, -18.92234920697049, 1.320258555276565, -18.92234920697049, 1.666878572882711, -18.92234920697049, 2.507041698583423, 1.24597646333707, -18.92234920697049, -18.92234920697049, -18.92234920697049)
)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($ci = 1; $ci -lt $entry.Length; $ci++) {
        $ws.Cells.Item($r, $ci + 1).Value = $entry[$ci]
    }
}
